# Add a second Q&A paragraph ("Q2. What do you understand by traversing a
# linked list? / Ans. ...") right before the final empty paragraph at the
# end of the document (i.e. between the two trailing empty paragraphs that
# sit just above the sectPr).
#
# The new paragraph needs three separate runs:
#   1) "Q2. "
#   2) "What do you understand by traversing a linked list?"
#   3) a line break followed by "Ans. Traverse means travel. ..."
#
# A straightforward Range.InsertAfter (or Selection.TypeText) sequence
# collapses adjacent same-format runs into a single run on save, so instead
# each segment is typed into its own temporary paragraph and the paragraph
# marks between them are then deleted -- merging the paragraphs back into
# one while keeping each segment as its own <w:r>.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$targetIndex = $count - 1   # the first of the two trailing empty paragraphs

$anchor = $d.Paragraphs.Item($targetIndex)
$anchor.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$insPos = $newPara.Range.Start

# --- Run 1: "Q2. " ---
$r1 = $d.Range($insPos, $insPos)
$r1.InsertAfter("Q2. ")

# --- Run 2: "What do you understand by traversing a linked list?" ---
$r1end = $d.Range($r1.End, $r1.End)
$r1end.InsertParagraphAfter()
$tmpPara1 = $d.Paragraphs.Item($targetIndex + 2)
$r2 = $d.Range($tmpPara1.Range.Start, $tmpPara1.Range.Start)
$r2.InsertAfter("What do you understand by traversing a linked list?")

# --- Run 3: line break + "Ans. Traverse means travel. ..." ---
$r2end = $d.Range($r2.End, $r2.End)
$r2end.InsertParagraphAfter()
$tmpPara2 = $d.Paragraphs.Item($targetIndex + 3)
$r3 = $d.Range($tmpPara2.Range.Start, $tmpPara2.Range.Start)
$r3.InsertAfter([char]11 + "Ans. Traverse means travel. So, in a linked list traverse means travelling across the linked list element by element. Traversing helps us to search, delete or print elements. ")

# Merge the three temporary paragraphs back into a single paragraph by
# deleting the paragraph-mark characters between them. This preserves each
# typed segment as its own run instead of collapsing them into one run.
$p = $d.Paragraphs.Item($targetIndex + 1)
$mark1 = $d.Range($p.Range.End - 1, $p.Range.End)
$mark1.Delete()

$p2 = $d.Paragraphs.Item($targetIndex + 1)
$mark2 = $d.Range($p2.Range.End - 1, $p2.Range.End)
$mark2.Delete()

Write-Output "Inserted Q2 paragraph; document now has $($d.Paragraphs.Count) paragraphs."
